$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data sorted descending by value, with the two lowest-value rows
# (Vietnamese, Uzbek) dropped entirely.
$data = @(
    @("English", 26.72015198895935),
    @("Spanish", 8.122719392654663),
    @("Chinese", 7.692558688093446),
    @("Japanese", 7.580531673219575),
    @("German", 6.034929331336038),
    @("Arabic", 4.599224186790686),
    @("Portuguese", 3.853843549431342),
    @("French", 3.529082219073772),
    @("Italian", 3.414148384827087),
    @("Russian", 3.168879429919302),
    @("Malay-Indonesian", 2.994983643764452),
    @("Dutch", 1.632381422061824),
    @("Korean", 1.508704004642031),
    @("Turkish", 1.448065486269834),
    @("Persian", 1.425608745673678),
    @("Thai", 1.003625060446002),
    @("Polish", 0.8692384302056134),
    @("Urdu", 0.7992758897879318),
    @("Swedish", 0.5053471948302737),
    @("Bengali", 0.4136378529949986)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# The table used to go down to row 23 (Uzbek, Vietnamese); remove the
# now-unused trailing rows 22 and 23.
$ws.Range("A22:B23").Delete() | Out-Null
